# Fix logging system configuration
# Append a new row (75) of log data to each of the four worksheets,
# mirroring the structure of the existing rows.

$wb = $excel.ActiveWorkbook

$dateValue = 45861.43658564815

$rows = @(
    @{ Sheet = "DE_LFT_#1"; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x44"; E = "0x14"; F = 380;  G = [double]"7.598631275147109e+23"; H = 324; I = 14 },
    @{ Sheet = "DE_LFT_#2"; B = "0x01,0x7c"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x44"; E = "0xe";  F = 380;  G = [double]"5.68432987514711e+23";  H = 324; I = 14 },
    @{ Sheet = "DE_PLT_#1"; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x79"; E = "0x7";  F = 130;  G = [double]"5.68631262647114e+23";  H = 121; I = 7 },
    @{ Sheet = "DE_PLT_#2"; B = "0x00,0x82"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x79"; E = "0x3";  F = 130;  G = [double]"9.85046333984776e+23";  H = 121; I = 3 }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    $ws.Range("A75").Value = $dateValue
    $ws.Range("A75").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B75").Value = $row.B
    $ws.Range("C75").Value = $row.C
    $ws.Range("D75").Value = $row.D
    $ws.Range("E75").Value = $row.E
    $ws.Range("F75").Value = $row.F
    $ws.Range("G75").Value = $row.G
    $ws.Range("H75").Value = $row.H
    $ws.Range("I75").Value = $row.I
}
